$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "id_electricity_consumption" header/value to "id_electricity"
$ws.Range("B1").Value = "id_electricity"

# Update the table column header to match
$table = $ws.ListObjects.Item(1)
$table.ListColumns.Item(2).Name = "id_electricity"

# Move active selection to B1
$ws.Range("B1").Select()
